$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.749.29'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '2.593.92'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '552.61'
$ws.Range('E5').Value = '  +2.76%  '
$ws.Range('D6').Value = '142.96'
$ws.Range('E6').Value = '  -2.26%  '
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').Value = '  +4.91%  '
$ws.Range('D9').Value = '6.77'
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('D10').Value = '0.100'
$ws.Range('E10').Value = '  -2.02%  '
$ws.Range('E11').Value = '  +5.06%  '
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('D13').Value = '3.055.10'
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('D14').Value = '58.726.21'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').Value = '20.82'
$ws.Range('E15').Value = '  -2.25%  '
$ws.Range('D16').Value = '2.609.19'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('D19').Value = '337.05'
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('E20').Value = '  -2.57%  '
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '66.78'
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').Value = '0.427'
$ws.Range('E24').Value = '  +2.41%  '
$ws.Range('D25').Value = '0.994'
$ws.Range('E25').Value = '  -0.52%  '
$ws.Range('E26').Value = '  -3.57%  '
$ws.Range('D27').Value = '7.12'
$ws.Range('D28').Value = '0.0₃0750'
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('D31').Value = '5.98'
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('D32').Value = '154.78'
$ws.Range('E32').Value = '  +2.64%  '
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('E34').Value = '  -2.18%  '
$ws.Range('D35').Value = '0.891'
$ws.Range('E35').Value = '  +6.53%  '
$ws.Range('E36').Value = '  -0.92%  '
$ws.Range('D37').Value = '36.85'
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').Value = '0.846'
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '1.46'
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('D40').Value = '3.59'
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('D41').Value = '283.39'
$ws.Range('E41').Value = '  -1.32%  '
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('E45').Value = '  -1.09%  '
$ws.Range('E46').Value = '  -1.47%  '
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = '1.943.72'
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('D49').Value = '118.10'
$ws.Range('E49').Value = '  +6.32%  '
$ws.Range('D50').Value = '17.85'
$ws.Range('E50').Value = '  -2.78%  '
$ws.Range('E51').Value = '  -3.40%  '
